$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 218 = verbatim copy of the old row 217 data (the whole 45-row weekly
# block shifts down by one row to make room for the new week inserted at row 173).
# Row 218: brand-new row (did not exist before)
$ws.Cells.Item(218, 1).Value = 5
$ws.Cells.Item(218, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(218, 3).Value = 'Maule'
$ws.Cells.Item(218, 4).Value = 44272
$ws.Cells.Item(218, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(218, 5).Value = 7
$ws.Cells.Item(218, 6).Value = 100112003
$ws.Cells.Item(218, 7).Value = 'Ajo'
$ws.Cells.Item(218, 8).Value = 'Chino'
$ws.Cells.Item(218, 9).Value = 'Primera'
$ws.Cells.Item(218, 10).Value = 200
$ws.Cells.Item(218, 11).Value = 12000
$ws.Cells.Item(218, 12).Value = 12000
$ws.Cells.Item(218, 13).Value = 12000
$ws.Cells.Item(218, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(218, 15).Value = 'China'
$ws.Cells.Item(218, 16).Value = 1200
$ws.Cells.Item(218, 17).Value = 10
$ws.Cells.Item(218, 18).Value = 'Hortaliza'

# Row 217
$ws.Cells.Item(217, 4).Value = 44518
$ws.Cells.Item(217, 11).Value = 20000
$ws.Cells.Item(217, 12).Value = 20000
$ws.Cells.Item(217, 13).Value = 20000
$ws.Cells.Item(217, 16).Value = 2000

# Row 216
$ws.Cells.Item(216, 4).Value = 44313
$ws.Cells.Item(216, 10).Value = 300
$ws.Cells.Item(216, 11).Value = 14000
$ws.Cells.Item(216, 12).Value = 14000
$ws.Cells.Item(216, 13).Value = 14000
$ws.Cells.Item(216, 16).Value = 1400

# Row 215
$ws.Cells.Item(215, 4).Value = 44194
$ws.Cells.Item(215, 11).Value = 9000
$ws.Cells.Item(215, 12).Value = 9000
$ws.Cells.Item(215, 13).Value = 9000
$ws.Cells.Item(215, 16).Value = 900

# Row 214
$ws.Cells.Item(214, 4).Value = 44371
$ws.Cells.Item(214, 11).Value = 13000
$ws.Cells.Item(214, 12).Value = 13000
$ws.Cells.Item(214, 13).Value = 13000
$ws.Cells.Item(214, 16).Value = 1300

# Row 213
$ws.Cells.Item(213, 4).Value = 44357
$ws.Cells.Item(213, 10).Value = 200

# Row 212
$ws.Cells.Item(212, 4).Value = 44463
$ws.Cells.Item(212, 11).Value = 15000
$ws.Cells.Item(212, 12).Value = 15000
$ws.Cells.Item(212, 13).Value = 15000
$ws.Cells.Item(212, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(212, 16).Value = 1500

# Row 211
$ws.Cells.Item(211, 4).Value = 44461

# Row 210
$ws.Cells.Item(210, 4).Value = 44511
$ws.Cells.Item(210, 11).Value = 20000
$ws.Cells.Item(210, 12).Value = 20000
$ws.Cells.Item(210, 13).Value = 20000
$ws.Cells.Item(210, 16).Value = 2000

# Row 209
$ws.Cells.Item(209, 4).Value = 44302
$ws.Cells.Item(209, 11).Value = 14000
$ws.Cells.Item(209, 12).Value = 14000
$ws.Cells.Item(209, 13).Value = 14000
$ws.Cells.Item(209, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(209, 16).Value = 1400

# Row 208
$ws.Cells.Item(208, 4).Value = 44536
$ws.Cells.Item(208, 11).Value = 18000
$ws.Cells.Item(208, 12).Value = 18000
$ws.Cells.Item(208, 13).Value = 18000
$ws.Cells.Item(208, 16).Value = 1800

# Row 207
$ws.Cells.Item(207, 4).Value = 44398
$ws.Cells.Item(207, 11).Value = 13000
$ws.Cells.Item(207, 12).Value = 13000
$ws.Cells.Item(207, 13).Value = 13000
$ws.Cells.Item(207, 16).Value = 1300

# Row 206
$ws.Cells.Item(206, 4).Value = 44473
$ws.Cells.Item(206, 11).Value = 16000
$ws.Cells.Item(206, 12).Value = 16000
$ws.Cells.Item(206, 13).Value = 16000
$ws.Cells.Item(206, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(206, 16).Value = 1600

# Row 205
$ws.Cells.Item(205, 4).Value = 44195
$ws.Cells.Item(205, 10).Value = 300
$ws.Cells.Item(205, 11).Value = 9000
$ws.Cells.Item(205, 12).Value = 9000
$ws.Cells.Item(205, 13).Value = 9000
$ws.Cells.Item(205, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(205, 16).Value = 900

# Row 204
$ws.Cells.Item(204, 4).Value = 44270
$ws.Cells.Item(204, 11).Value = 12000
$ws.Cells.Item(204, 12).Value = 12000
$ws.Cells.Item(204, 13).Value = 12000
$ws.Cells.Item(204, 15).Value = 'Región del Maule'
$ws.Cells.Item(204, 16).Value = 1200

# Row 203
$ws.Cells.Item(203, 4).Value = 44452
$ws.Cells.Item(203, 11).Value = 15000
$ws.Cells.Item(203, 12).Value = 15000
$ws.Cells.Item(203, 13).Value = 15000
$ws.Cells.Item(203, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(203, 15).Value = 'China'
$ws.Cells.Item(203, 16).Value = 1500

# Row 202
$ws.Cells.Item(202, 4).Value = 44329

# Row 201
$ws.Cells.Item(201, 4).Value = 44224
$ws.Cells.Item(201, 10).Value = 200
$ws.Cells.Item(201, 11).Value = 13000
$ws.Cells.Item(201, 12).Value = 13000
$ws.Cells.Item(201, 13).Value = 13000
$ws.Cells.Item(201, 16).Value = 1300

# Row 200
$ws.Cells.Item(200, 4).Value = 44340
$ws.Cells.Item(200, 14).Value = '$/caja 10 kilos'

# Row 199
$ws.Cells.Item(199, 4).Value = 44428
$ws.Cells.Item(199, 11).Value = 14000
$ws.Cells.Item(199, 12).Value = 14000
$ws.Cells.Item(199, 13).Value = 14000
$ws.Cells.Item(199, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(199, 16).Value = 1400

# Row 198
$ws.Cells.Item(198, 4).Value = 44432

# Row 197
$ws.Cells.Item(197, 4).Value = 44298
$ws.Cells.Item(197, 10).Value = 300
$ws.Cells.Item(197, 14).Value = '$/caja 10 kilos'

# Row 196
$ws.Cells.Item(196, 4).Value = 44258
$ws.Cells.Item(196, 10).Value = 200
$ws.Cells.Item(196, 11).Value = 13000
$ws.Cells.Item(196, 12).Value = 13000
$ws.Cells.Item(196, 13).Value = 13000
$ws.Cells.Item(196, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(196, 16).Value = 1300

# Row 195
$ws.Cells.Item(195, 4).Value = 44459
$ws.Cells.Item(195, 10).Value = 300
$ws.Cells.Item(195, 11).Value = 15000
$ws.Cells.Item(195, 12).Value = 15000
$ws.Cells.Item(195, 13).Value = 15000
$ws.Cells.Item(195, 15).Value = 'Región del Maule'
$ws.Cells.Item(195, 16).Value = 1500

# Row 194
$ws.Cells.Item(194, 4).Value = 44411
$ws.Cells.Item(194, 11).Value = 13000
$ws.Cells.Item(194, 12).Value = 13000
$ws.Cells.Item(194, 13).Value = 13000
$ws.Cells.Item(194, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(194, 15).Value = 'China'
$ws.Cells.Item(194, 16).Value = 1300

# Row 193
$ws.Cells.Item(193, 4).Value = 44257

# Row 192
$ws.Cells.Item(192, 4).Value = 44238
$ws.Cells.Item(192, 11).Value = 12000
$ws.Cells.Item(192, 12).Value = 12000
$ws.Cells.Item(192, 13).Value = 12000
$ws.Cells.Item(192, 15).Value = 'Región del Maule'
$ws.Cells.Item(192, 16).Value = 1200

# Row 191
$ws.Cells.Item(191, 4).Value = 44413
$ws.Cells.Item(191, 11).Value = 14000
$ws.Cells.Item(191, 12).Value = 14000
$ws.Cells.Item(191, 13).Value = 14000
$ws.Cells.Item(191, 15).Value = 'China'
$ws.Cells.Item(191, 16).Value = 1400

# Row 190
$ws.Cells.Item(190, 4).Value = 44405
$ws.Cells.Item(190, 11).Value = 13000
$ws.Cells.Item(190, 12).Value = 13000
$ws.Cells.Item(190, 13).Value = 13000
$ws.Cells.Item(190, 16).Value = 1300

# Row 189
$ws.Cells.Item(189, 4).Value = 44529
$ws.Cells.Item(189, 10).Value = 200
$ws.Cells.Item(189, 11).Value = 22000
$ws.Cells.Item(189, 12).Value = 22000
$ws.Cells.Item(189, 13).Value = 22000
$ws.Cells.Item(189, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(189, 16).Value = 2200

# Row 188
$ws.Cells.Item(188, 4).Value = 44350
$ws.Cells.Item(188, 11).Value = 13000
$ws.Cells.Item(188, 12).Value = 13000
$ws.Cells.Item(188, 13).Value = 13000
$ws.Cells.Item(188, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(188, 16).Value = 1300

# Row 187
$ws.Cells.Item(187, 4).Value = 44446
$ws.Cells.Item(187, 11).Value = 14000
$ws.Cells.Item(187, 12).Value = 14000
$ws.Cells.Item(187, 13).Value = 14000
$ws.Cells.Item(187, 16).Value = 1400

# Row 186
$ws.Cells.Item(186, 4).Value = 44474
$ws.Cells.Item(186, 10).Value = 300
$ws.Cells.Item(186, 11).Value = 16000
$ws.Cells.Item(186, 12).Value = 16000
$ws.Cells.Item(186, 13).Value = 16000
$ws.Cells.Item(186, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(186, 16).Value = 1600

# Row 185
$ws.Cells.Item(185, 4).Value = 44187
$ws.Cells.Item(185, 10).Value = 200
$ws.Cells.Item(185, 11).Value = 8000
$ws.Cells.Item(185, 12).Value = 8000
$ws.Cells.Item(185, 13).Value = 8000
$ws.Cells.Item(185, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(185, 16).Value = 800

# Row 184
$ws.Cells.Item(184, 4).Value = 44281
$ws.Cells.Item(184, 11).Value = 13000
$ws.Cells.Item(184, 12).Value = 13000
$ws.Cells.Item(184, 13).Value = 13000
$ws.Cells.Item(184, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(184, 16).Value = 1300

# Row 183
$ws.Cells.Item(183, 4).Value = 44308
$ws.Cells.Item(183, 10).Value = 300
$ws.Cells.Item(183, 11).Value = 14000
$ws.Cells.Item(183, 12).Value = 14000
$ws.Cells.Item(183, 13).Value = 14000
$ws.Cells.Item(183, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(183, 16).Value = 1400

# Row 182
$ws.Cells.Item(182, 4).Value = 44278
$ws.Cells.Item(182, 11).Value = 13000
$ws.Cells.Item(182, 12).Value = 13000
$ws.Cells.Item(182, 13).Value = 13000
$ws.Cells.Item(182, 16).Value = 1300

# Row 181
$ws.Cells.Item(181, 4).Value = 44385
$ws.Cells.Item(181, 10).Value = 500

# Row 180
$ws.Cells.Item(180, 4).Value = 44321
$ws.Cells.Item(180, 10).Value = 300
$ws.Cells.Item(180, 11).Value = 15000
$ws.Cells.Item(180, 12).Value = 15000
$ws.Cells.Item(180, 13).Value = 15000
$ws.Cells.Item(180, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(180, 16).Value = 1500

# Row 179
$ws.Cells.Item(179, 4).Value = 44543
$ws.Cells.Item(179, 10).Value = 200
$ws.Cells.Item(179, 11).Value = 20000
$ws.Cells.Item(179, 12).Value = 20000
$ws.Cells.Item(179, 13).Value = 20000
$ws.Cells.Item(179, 16).Value = 2000

# Row 178
$ws.Cells.Item(178, 11).Value = 18000
$ws.Cells.Item(178, 12).Value = 18000
$ws.Cells.Item(178, 13).Value = 18000
$ws.Cells.Item(178, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(178, 16).Value = 1800

# Row 177
$ws.Cells.Item(177, 4).Value = 44414
$ws.Cells.Item(177, 10).Value = 300
$ws.Cells.Item(177, 11).Value = 14000
$ws.Cells.Item(177, 12).Value = 14000
$ws.Cells.Item(177, 13).Value = 14000
$ws.Cells.Item(177, 16).Value = 1400

# Row 176
$ws.Cells.Item(176, 4).Value = 44420

# Row 175
$ws.Cells.Item(175, 4).Value = 44519
$ws.Cells.Item(175, 10).Value = 200
$ws.Cells.Item(175, 11).Value = 21000
$ws.Cells.Item(175, 12).Value = 21000
$ws.Cells.Item(175, 13).Value = 21000
$ws.Cells.Item(175, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(175, 16).Value = 2100

# Row 174
$ws.Cells.Item(174, 4).Value = 44386
$ws.Cells.Item(174, 10).Value = 800
$ws.Cells.Item(174, 11).Value = 13000
$ws.Cells.Item(174, 12).Value = 13000
$ws.Cells.Item(174, 13).Value = 13000
$ws.Cells.Item(174, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(174, 16).Value = 1300

# Row 173
$ws.Cells.Item(173, 4).Value = 44551
$ws.Cells.Item(173, 10).Value = 200
$ws.Cells.Item(173, 11).Value = 20000
$ws.Cells.Item(173, 12).Value = 20000
$ws.Cells.Item(173, 13).Value = 20000
$ws.Cells.Item(173, 16).Value = 2000
